# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) values for the first data row
# on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-12 10:47:20"
$wsZhCn.Range("H2").Value = "2016-03-12 10:47:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-12 10:47:23"
$wsDeDe.Range("H2").Value = "2016-03-12 10:47:42"
